# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" for the files that were just
# (re)handed off, and flag them with the "ht" (handoff type) priority
# marker on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 12, 13)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-14 16:29:10"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-14 16:28:58"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-14 16:29:10"
}
